# Add a new "Search" worksheet as the last tab, and populate it with the
# new search-related test case data ("search test cases added").

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Search"

# Column C (category_value) is text-formatted; set that before writing any
# values so the numeric-looking entries ("27", "18") are stored as text
# rather than numbers.
$ws.Range("C2:C10").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "action"
$ws.Range("B1").Value = "searched_text"
$ws.Range("C1").Value = "category_value"

# Row 2
$ws.Range("A2").Value = "searched_product_text_validation"
$ws.Range("B2").Value = "iphone"

# Row 3
$ws.Range("A3").Value = "no_product_found"
$ws.Range("B3").Value = "productnotexists"

# Row 4 - empty search text, entered as a lone apostrophe (forces an empty
# text cell with the quote-prefix style, same as typing `'` into Excel).
$ws.Range("A4").Value = "no_product_found"
$ws.Range("B4").Value = "'"

# Row 5
$ws.Range("A5").Value = "searched_product_after_login"
$ws.Range("B5").Value = "samsung"

# Row 6
$ws.Range("A6").Value = "multiple_product_found"
$ws.Range("B6").Value = "e"

# Row 7
$ws.Range("A7").Value = "multiple_product_found"
$ws.Range("B7").Value = "samsung"

# Row 8
$ws.Range("A8").Value = "search_product_under_search_criteria"
$ws.Range("B8").Value = "iMac"

# Row 9
$ws.Range("A9").Value = "search_product_by_category_product_found"
$ws.Range("B9").Value = "imac"
$ws.Range("C9").Value = "27"

# Row 10
$ws.Range("A10").Value = "search_product_by_category_no_product_found"
$ws.Range("B10").Value = "imac"
$ws.Range("C10").Value = "18"

# Row 11 - C11 stays a real number (no text format applied to this cell).
$ws.Range("A11").Value = "search_product_check_sub_category"
$ws.Range("B11").Value = "imac"
$ws.Range("C11").Value = 20

# Row 12
$ws.Range("A12").Value = "product_compare_navigation_check"
$ws.Range("B12").Value = "HTC Touch HD"

$ws.Columns.Item(1).AutoFit()
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(3).AutoFit()

# Selection / active-sheet bookkeeping to mirror the authored edit: the new
# Search sheet becomes the active tab with B12 selected, while Registration
# keeps B2 as its (non-active) selection.
$registration = $wb.Worksheets.Item("Registration")
[void]$registration.Range("B2").Select()

$ws.Activate()
[void]$ws.Range("B12").Select()
